$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "trainingimages/08_tipako"
$ws.Range("B2").Value = "pngimages/08_bell.png"
$ws.Range("C2").Value = "trainingimages/20_tatito"
$ws.Range("D2").Value = "pngimages/20_pizza.png"

$ws.Range("A3").Value = "trainingimages/06_titoka"
$ws.Range("B3").Value = "pngimages/06_tent.png"
$ws.Range("C3").Value = "trainingimages/04_kitoti"
$ws.Range("D3").Value = "pngimages/04_ladder.png"

$ws.Range("A4").Value = "trainingimages/26_kapako"
$ws.Range("B4").Value = "pngimages/26_pineapple.png"
$ws.Range("C4").Value = "trainingimages/12_pokika"
$ws.Range("D4").Value = "pngimages/12_pie.png"
